$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.485.45"
$ws.Range("E2").Value = "  -0.83%  "

$ws.Range("D3").Value = "1.622.57"
$ws.Range("E3").Value = "  -0.02%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.55"
$ws.Range("E5").Value = "  -0.47%  "

$ws.Range("E6").Value = "  -0.76%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E9").Value = "  -0.45%  "

$ws.Range("E10").Value = "  -0.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0855"
$ws.Range("E11").Value = "  -0.41%  "

$ws.Range("D12").Value = "1.850.65"
$ws.Range("E12").Value = "  -0.04%  "

$ws.Range("D13").Value = "1.625.39"
$ws.Range("E13").Value = "  +0.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.04"
$ws.Range("E14").Value = "  -0.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.510"
$ws.Range("E15").Value = "  -0.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.97"
$ws.Range("E16").Value = "  -1.69%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "26.494.87"
$ws.Range("E17").Value = "  -0.80%  "

$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "233.76"
$ws.Range("E18").Value = "  +1.32%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.72"
$ws.Range("E19").Value = "  -0.23%  "

$ws.Range("D20").Value = "0.0₃0725"
$ws.Range("E20").Value = "  -0.51%  "

$ws.Range("E21").Value = "  -0.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.31"
$ws.Range("E22").Value = "  -1.81%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.22"
$ws.Range("E23").Value = "  -0.58%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.14"
$ws.Range("E24").Value = "  +0.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.04"
$ws.Range("E25").Value = "  +0.03%  "

$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.06"
$ws.Range("E27").Value = "  -0.07%  "

$ws.Range("E28").Value = "  -0.61%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.60"
$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0494"
$ws.Range("E30").Value = "  -0.92%  "

$ws.Range("E31").Value = "  -0.40%  "

$ws.Range("D32").Value = "1.518.59"
$ws.Range("E32").Value = "  +3.67%  "

$ws.Range("E33").Value = "  +0.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.99"
$ws.Range("E34").Value = "  -0.44%  "

$ws.Range("E35").Value = "  +1.99%  "

$ws.Range("E36").Value = "  -0.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.567"
$ws.Range("E37").Value = "  -1.12%  "

$ws.Range("E38").Value = "  -0.62%  "

$ws.Range("E39").Value = "  -0.66%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.85"
$ws.Range("E40").Value = "  -1.74%  "

$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.22"
$ws.Range("E42").Value = "  +0.51%  "

$ws.Range("D43").Value = "1.761.77"
$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.57"
$ws.Range("E44").Value = "  +0.77%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.760"
$ws.Range("E45").Value = "  -0.47%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.905"
$ws.Range("E46").Value = "  -5.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.72"
$ws.Range("E47").Value = "  +1.82%  "

$ws.Range("E48").Value = "  +0.51%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0501"
$ws.Range("E49").Value = "  -0.85%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.55"
$ws.Range("E50").Value = "  +0.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0962"
$ws.Range("E51").Value = "  -0.09%  "
